$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 155,4
$data[0,0] = 'Programme ID'; $data[0,1] = 'Module ID'; $data[0,2] = 'Start Year'; $data[0,3] = 'End Year'
$data[1,0] = 'BIT'; $data[1,1] = 'BIS2104'; $data[1,2] = 2022; $data[1,3] = $null
$data[2,0] = 'BIT'; $data[2,1] = 'BIS2216'; $data[2,2] = 2022; $data[2,3] = $null
$data[3,0] = 'BIT'; $data[3,1] = 'BIS3106'; $data[3,2] = 2022; $data[3,3] = $null
$data[4,0] = 'BIT'; $data[4,1] = 'CSC1024'; $data[4,2] = 2022; $data[4,3] = $null
$data[5,0] = 'BIT'; $data[5,1] = 'CSC1202'; $data[5,2] = 2022; $data[5,3] = $null
$data[6,0] = 'BIT'; $data[6,1] = 'CSC2103'; $data[6,2] = 2022; $data[6,3] = $null
$data[7,0] = 'BIT'; $data[7,1] = 'CSC2104'; $data[7,2] = 2022; $data[7,3] = $null
$data[8,0] = 'BIT'; $data[8,1] = 'CSC3024'; $data[8,2] = 2022; $data[8,3] = $null
$data[9,0] = 'BIT'; $data[9,1] = 'CSC3044'; $data[9,2] = 2022; $data[9,3] = $null
$data[10,0] = 'BIT'; $data[10,1] = 'CSC3064'; $data[10,2] = 2022; $data[10,3] = $null
$data[11,0] = 'BIT'; $data[11,1] = 'CSC3206'; $data[11,2] = 2022; $data[11,3] = $null
$data[12,0] = 'BIT'; $data[12,1] = 'CSC3209'; $data[12,2] = 2022; $data[12,3] = $null
$data[13,0] = 'BIT'; $data[13,1] = 'ENG1044'; $data[13,2] = 2022; $data[13,3] = $null
$data[14,0] = 'BIT'; $data[14,1] = 'ETP2014'; $data[14,2] = 2022; $data[14,3] = $null
$data[15,0] = 'BIT'; $data[15,1] = 'MTH1114'; $data[15,2] = 2022; $data[15,3] = $null
$data[16,0] = 'BIT'; $data[16,1] = 'MU4 2422'; $data[16,2] = 2022; $data[16,3] = $null
$data[17,0] = 'BIT'; $data[17,1] = 'NET1014'; $data[17,2] = 2022; $data[17,3] = $null
$data[18,0] = 'BIT'; $data[18,1] = 'NET2102'; $data[18,2] = 2022; $data[18,3] = $null
$data[19,0] = 'BIT'; $data[19,1] = 'NET2201'; $data[19,2] = 2022; $data[19,3] = $null
$data[20,0] = 'BIT'; $data[20,1] = 'NET3106'; $data[20,2] = 2022; $data[20,3] = $null
$data[21,0] = 'BIT'; $data[21,1] = 'NET3204'; $data[21,2] = 2022; $data[21,3] = $null
$data[22,0] = 'BIT'; $data[22,1] = 'NET3207'; $data[22,2] = 2022; $data[22,3] = $null
$data[23,0] = 'BIT'; $data[23,1] = 'PRG1203'; $data[23,2] = 2022; $data[23,3] = $null
$data[24,0] = 'BIT'; $data[24,1] = 'PRG2104'; $data[24,2] = 2022; $data[24,3] = $null
$data[25,0] = 'BIT'; $data[25,1] = 'PRG2224'; $data[25,2] = 2022; $data[25,3] = $null
$data[26,0] = 'BIT'; $data[26,1] = 'PRJ3213'; $data[26,2] = 2022; $data[26,3] = $null
$data[27,0] = 'BIT'; $data[27,1] = 'PRJ3223'; $data[27,2] = 2022; $data[27,3] = $null
$data[28,0] = 'BIT'; $data[28,1] = 'SEC3014'; $data[28,2] = 2022; $data[28,3] = $null
$data[29,0] = 'BIT'; $data[29,1] = 'SEG1201'; $data[29,2] = 2022; $data[29,3] = $null
$data[30,0] = 'BIT'; $data[30,1] = 'SEG2102'; $data[30,2] = 2022; $data[30,3] = $null
$data[31,0] = 'BIT'; $data[31,1] = 'SEG2202'; $data[31,2] = 2022; $data[31,3] = $null
$data[32,0] = 'BIT'; $data[32,1] = 'SEG3203'; $data[32,2] = 2022; $data[32,3] = $null
$data[33,0] = 'BIT'; $data[33,1] = 'WEB1201'; $data[33,2] = 2022; $data[33,3] = $null
$data[34,0] = 'BIT'; $data[34,1] = 'WEB2202'; $data[34,2] = 2022; $data[34,3] = $null
$data[35,0] = 'BCS'; $data[35,1] = 'BIS2216'; $data[35,2] = 2022; $data[35,3] = $null
$data[36,0] = 'BCS'; $data[36,1] = 'CSC1024'; $data[36,2] = 2022; $data[36,3] = $null
$data[37,0] = 'BCS'; $data[37,1] = 'CSC1202'; $data[37,2] = 2022; $data[37,3] = $null
$data[38,0] = 'BCS'; $data[38,1] = 'CSC2014'; $data[38,2] = 2022; $data[38,3] = $null
$data[39,0] = 'BCS'; $data[39,1] = 'CSC2103'; $data[39,2] = 2022; $data[39,3] = $null
$data[40,0] = 'BCS'; $data[40,1] = 'CSC2104'; $data[40,2] = 2022; $data[40,3] = $null
$data[41,0] = 'BCS'; $data[41,1] = 'CSC3024'; $data[41,2] = 2022; $data[41,3] = $null
$data[42,0] = 'BCS'; $data[42,1] = 'CSC3034'; $data[42,2] = 2022; $data[42,3] = $null
$data[43,0] = 'BCS'; $data[43,1] = 'CSC3044'; $data[43,2] = 2022; $data[43,3] = $null
$data[44,0] = 'BCS'; $data[44,1] = 'CSC3064'; $data[44,2] = 2022; $data[44,3] = $null
$data[45,0] = 'BCS'; $data[45,1] = 'CSC3206'; $data[45,2] = 2022; $data[45,3] = $null
$data[46,0] = 'BCS'; $data[46,1] = 'CSC3209'; $data[46,2] = 2022; $data[46,3] = $null
$data[47,0] = 'BCS'; $data[47,1] = 'ENG1044'; $data[47,2] = 2022; $data[47,3] = $null
$data[48,0] = 'BCS'; $data[48,1] = 'ETP2014'; $data[48,2] = 2022; $data[48,3] = $null
$data[49,0] = 'BCS'; $data[49,1] = 'IST2334'; $data[49,2] = 2022; $data[49,3] = $null
$data[50,0] = 'BCS'; $data[50,1] = 'MAT1013'; $data[50,2] = 2022; $data[50,3] = $null
$data[51,0] = 'BCS'; $data[51,1] = 'MTH1114'; $data[51,2] = 2022; $data[51,3] = $null
$data[52,0] = 'BCS'; $data[52,1] = 'MU4 2422'; $data[52,2] = 2022; $data[52,3] = $null
$data[53,0] = 'BCS'; $data[53,1] = 'NET1014'; $data[53,2] = 2022; $data[53,3] = $null
$data[54,0] = 'BCS'; $data[54,1] = 'NET2102'; $data[54,2] = 2022; $data[54,3] = $null
$data[55,0] = 'BCS'; $data[55,1] = 'NET2201'; $data[55,2] = 2022; $data[55,3] = $null
$data[56,0] = 'BCS'; $data[56,1] = 'NET3204'; $data[56,2] = 2022; $data[56,3] = $null
$data[57,0] = 'BCS'; $data[57,1] = 'PRG1203'; $data[57,2] = 2022; $data[57,3] = $null
$data[58,0] = 'BCS'; $data[58,1] = 'PRG2104'; $data[58,2] = 2022; $data[58,3] = $null
$data[59,0] = 'BCS'; $data[59,1] = 'PRG2205'; $data[59,2] = 2022; $data[59,3] = $null
$data[60,0] = 'BCS'; $data[60,1] = 'PRG2214'; $data[60,2] = 2022; $data[60,3] = $null
$data[61,0] = 'BCS'; $data[61,1] = 'PRJ3213'; $data[61,2] = 2022; $data[61,3] = $null
$data[62,0] = 'BCS'; $data[62,1] = 'PRJ3223'; $data[62,2] = 2022; $data[62,3] = $null
$data[63,0] = 'BCS'; $data[63,1] = 'SEG1201'; $data[63,2] = 2022; $data[63,3] = $null
$data[64,0] = 'BCS'; $data[64,1] = 'SEG2102'; $data[64,2] = 2022; $data[64,3] = $null
$data[65,0] = 'BCS'; $data[65,1] = 'SEG2202'; $data[65,2] = 2022; $data[65,3] = $null
$data[66,0] = 'BCS'; $data[66,1] = 'SEG3203'; $data[66,2] = 2022; $data[66,3] = $null
$data[67,0] = 'BCS'; $data[67,1] = 'WEB1201'; $data[67,2] = 2022; $data[67,3] = $null
$data[68,0] = 'BSDA'; $data[68,1] = 'BIS1014'; $data[68,2] = 2022; $data[68,3] = $null
$data[69,0] = 'BSDA'; $data[69,1] = 'BIS3106'; $data[69,2] = 2022; $data[69,3] = $null
$data[70,0] = 'BSDA'; $data[70,1] = 'BIS3216'; $data[70,2] = 2022; $data[70,3] = $null
$data[71,0] = 'BSDA'; $data[71,1] = 'BIS3218'; $data[71,2] = 2022; $data[71,3] = $null
$data[72,0] = 'BSDA'; $data[72,1] = 'CSC1024'; $data[72,2] = 2022; $data[72,3] = $null
$data[73,0] = 'BSDA'; $data[73,1] = 'CSC1202'; $data[73,2] = 2022; $data[73,3] = $null
$data[74,0] = 'BSDA'; $data[74,1] = 'EAC2014'; $data[74,2] = 2022; $data[74,3] = $null
$data[75,0] = 'BSDA'; $data[75,1] = 'ENG1044'; $data[75,2] = 2022; $data[75,3] = $null
$data[76,0] = 'BSDA'; $data[76,1] = 'ETP2014'; $data[76,2] = 2022; $data[76,3] = $null
$data[77,0] = 'BSDA'; $data[77,1] = 'IST1014'; $data[77,2] = 2022; $data[77,3] = 2023
$data[78,0] = 'BSDA'; $data[78,1] = 'IST2034'; $data[78,2] = 2022; $data[78,3] = $null
$data[79,0] = 'BSDA'; $data[79,1] = 'IST2334'; $data[79,2] = 2022; $data[79,3] = 2023
$data[80,0] = 'BSDA'; $data[80,1] = 'MAN3154'; $data[80,2] = 2022; $data[80,3] = 2023
$data[81,0] = 'BSDA'; $data[81,1] = 'MAT1013'; $data[81,2] = 2024; $data[81,3] = $null
$data[82,0] = 'BSDA'; $data[82,1] = 'MU4 2422'; $data[82,2] = 2022; $data[82,3] = $null
$data[83,0] = 'BSDA'; $data[83,1] = 'NET1014'; $data[83,2] = 2022; $data[83,3] = $null
$data[84,0] = 'BSDA'; $data[84,1] = 'PRJ3213'; $data[84,2] = 2022; $data[84,3] = $null
$data[85,0] = 'BSDA'; $data[85,1] = 'PRJ3223'; $data[85,2] = 2022; $data[85,3] = $null
$data[86,0] = 'BSDA'; $data[86,1] = 'SEG1201'; $data[86,2] = 2022; $data[86,3] = $null
$data[87,0] = 'BSDA'; $data[87,1] = 'SEG3203'; $data[87,2] = 2022; $data[87,3] = $null
$data[88,0] = 'BSE'; $data[88,1] = 'BIS3106'; $data[88,2] = 2024; $data[88,3] = $null
$data[89,0] = 'BSE'; $data[89,1] = 'CSC1202'; $data[89,2] = 2022; $data[89,3] = $null
$data[90,0] = 'BSE'; $data[90,1] = 'CSC1024'; $data[90,2] = 2022; $data[90,3] = $null
$data[91,0] = 'BSE'; $data[91,1] = 'CSC2103'; $data[91,2] = 2022; $data[91,3] = $null
$data[92,0] = 'BSE'; $data[92,1] = 'CSC2104'; $data[92,2] = 2022; $data[92,3] = 2023
$data[93,0] = 'BSE'; $data[93,1] = 'OSS1014'; $data[93,2] = 2024; $data[93,3] = $null
$data[94,0] = 'BSE'; $data[94,1] = 'CSC3024'; $data[94,2] = 2022; $data[94,3] = $null
$data[95,0] = 'BSE'; $data[95,1] = 'CSC3044'; $data[95,2] = 2022; $data[95,3] = $null
$data[96,0] = 'BSE'; $data[96,1] = 'CSC3064'; $data[96,2] = 2022; $data[96,3] = $null
$data[97,0] = 'BSE'; $data[97,1] = 'CSC3206'; $data[97,2] = 2022; $data[97,3] = $null
$data[98,0] = 'BSE'; $data[98,1] = 'CSC3209'; $data[98,2] = 2022; $data[98,3] = $null
$data[99,0] = 'BSE'; $data[99,1] = 'ENG1044'; $data[99,2] = 2022; $data[99,3] = $null
$data[100,0] = 'BSE'; $data[100,1] = 'ETP2014'; $data[100,2] = 2022; $data[100,3] = $null
$data[101,0] = 'BSE'; $data[101,1] = 'IST1024'; $data[101,2] = 2022; $data[101,3] = 2023
$data[102,0] = 'BSE'; $data[102,1] = 'MAT1013'; $data[102,2] = 2023; $data[102,3] = $null
$data[103,0] = 'BSE'; $data[103,1] = 'MTH1114'; $data[103,2] = 2022; $data[103,3] = $null
$data[104,0] = 'BSE'; $data[104,1] = 'MU4 2422'; $data[104,2] = 2022; $data[104,3] = $null
$data[105,0] = 'BSE'; $data[105,1] = 'NET1014'; $data[105,2] = 2022; $data[105,3] = $null
$data[106,0] = 'BSE'; $data[106,1] = 'NET2102'; $data[106,2] = 2024; $data[106,3] = $null
$data[107,0] = 'BSE'; $data[107,1] = 'NET2201'; $data[107,2] = 2022; $data[107,3] = $null
$data[108,0] = 'BSE'; $data[108,1] = 'NET3106'; $data[108,2] = 2024; $data[108,3] = $null
$data[109,0] = 'BSE'; $data[109,1] = 'NET3204'; $data[109,2] = 2022; $data[109,3] = 2023
$data[110,0] = 'BSE'; $data[110,1] = 'PRG1203'; $data[110,2] = 2022; $data[110,3] = $null
$data[111,0] = 'BSE'; $data[111,1] = 'PRG2104'; $data[111,2] = 2022; $data[111,3] = $null
$data[112,0] = 'BSE'; $data[112,1] = 'PRJ3213'; $data[112,2] = 2022; $data[112,3] = $null
$data[113,0] = 'BSE'; $data[113,1] = 'PRJ3223'; $data[113,2] = 2022; $data[113,3] = $null
$data[114,0] = 'BSE'; $data[114,1] = 'SEC3014'; $data[114,2] = 2024; $data[114,3] = $null
$data[115,0] = 'BSE'; $data[115,1] = 'SEG1201'; $data[115,2] = 2022; $data[115,3] = $null
$data[116,0] = 'BSE'; $data[116,1] = 'SEG2102'; $data[116,2] = 2022; $data[116,3] = $null
$data[117,0] = 'BSE'; $data[117,1] = 'SEG2202'; $data[117,2] = 2022; $data[117,3] = $null
$data[118,0] = 'BSE'; $data[118,1] = 'SEG3203'; $data[118,2] = 2022; $data[118,3] = $null
$data[119,0] = 'BSE'; $data[119,1] = 'SWE2033'; $data[119,2] = 2022; $data[119,3] = $null
$data[120,0] = 'BSE'; $data[120,1] = 'SWE3024'; $data[120,2] = 2022; $data[120,3] = $null
$data[121,0] = 'BSE'; $data[121,1] = 'SWE3043'; $data[121,2] = 2022; $data[121,3] = $null
$data[122,0] = 'BSE'; $data[122,1] = 'SWE3053'; $data[122,2] = 2022; $data[122,3] = $null
$data[123,0] = 'BSE'; $data[123,1] = 'WEB1201'; $data[123,2] = 2022; $data[123,3] = $null
$data[124,0] = 'BSE'; $data[124,1] = 'WEB2202'; $data[124,2] = 2022; $data[124,3] = 2023
$data[125,0] = 'BCNS'; $data[125,1] = 'CSC1024'; $data[125,2] = 2022; $data[125,3] = $null
$data[126,0] = 'BCNS'; $data[126,1] = 'CSC1202'; $data[126,2] = 2022; $data[126,3] = $null
$data[127,0] = 'BCNS'; $data[127,1] = 'CSC2103'; $data[127,2] = 2022; $data[127,3] = $null
$data[128,0] = 'BCNS'; $data[128,1] = 'CSC2104'; $data[128,2] = 2022; $data[128,3] = 2023
$data[129,0] = 'BCNS'; $data[129,1] = 'OSS1014'; $data[129,2] = 2024; $data[129,3] = $null
$data[130,0] = 'BCNS'; $data[130,1] = 'CSC3024'; $data[130,2] = 2024; $data[130,3] = $null
$data[131,0] = 'BCNS'; $data[131,1] = 'CSC3044'; $data[131,2] = 2022; $data[131,3] = 2023
$data[132,0] = 'BCNS'; $data[132,1] = 'ENG1044'; $data[132,2] = 2022; $data[132,3] = $null
$data[133,0] = 'BCNS'; $data[133,1] = 'ETP2014'; $data[133,2] = 2022; $data[133,3] = $null
$data[134,0] = 'BCNS'; $data[134,1] = 'IST1014'; $data[134,2] = 2022; $data[134,3] = 2023
$data[135,0] = 'BCNS'; $data[135,1] = 'MAT1013'; $data[135,2] = 2024; $data[135,3] = $null
$data[136,0] = 'BCNS'; $data[136,1] = 'MMD3105'; $data[136,2] = 2022; $data[136,3] = 2023
$data[137,0] = 'BCNS'; $data[137,1] = 'MTH1114'; $data[137,2] = 2022; $data[137,3] = $null
$data[138,0] = 'BCNS'; $data[138,1] = 'MTH2103'; $data[138,2] = 2022; $data[138,3] = $null
$data[139,0] = 'BCNS'; $data[139,1] = 'MU4 2422'; $data[139,2] = 2022; $data[139,3] = $null
$data[140,0] = 'BCNS'; $data[140,1] = 'NET1014'; $data[140,2] = 2022; $data[140,3] = $null
$data[141,0] = 'BCNS'; $data[141,1] = 'NET2102'; $data[141,2] = 2022; $data[141,3] = $null
$data[142,0] = 'BCNS'; $data[142,1] = 'NET2201'; $data[142,2] = 2022; $data[142,3] = $null
$data[143,0] = 'BCNS'; $data[143,1] = 'NET3106'; $data[143,2] = 2022; $data[143,3] = $null
$data[144,0] = 'BCNS'; $data[144,1] = 'NET3204'; $data[144,2] = 2022; $data[144,3] = $null
$data[145,0] = 'BCNS'; $data[145,1] = 'NET3207'; $data[145,2] = 2022; $data[145,3] = 2023
$data[146,0] = 'BCNS'; $data[146,1] = 'PRG1203'; $data[146,2] = 2022; $data[146,3] = $null
$data[147,0] = 'BCNS'; $data[147,1] = 'PRJ3213'; $data[147,2] = 2022; $data[147,3] = $null
$data[148,0] = 'BCNS'; $data[148,1] = 'PRJ3223'; $data[148,2] = 2022; $data[148,3] = $null
$data[149,0] = 'BCNS'; $data[149,1] = 'SEC3014'; $data[149,2] = 2022; $data[149,3] = 2023
$data[150,0] = 'BCNS'; $data[150,1] = 'SEC3034'; $data[150,2] = 2022; $data[150,3] = $null
$data[151,0] = 'BCNS'; $data[151,1] = 'SEC3044'; $data[151,2] = 2022; $data[151,3] = $null
$data[152,0] = 'BCNS'; $data[152,1] = 'SEG1201'; $data[152,2] = 2022; $data[152,3] = $null
$data[153,0] = 'BCNS'; $data[153,1] = 'SEG3203'; $data[153,2] = 2022; $data[153,3] = $null
$data[154,0] = 'BCNS'; $data[154,1] = 'WEB1201'; $data[154,2] = 2022; $data[154,3] = $null

$ws.Range("A1:D155").Value = $data

